$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first point of the temperature profile is always room temperature,
# so the explicit (0.0 min, 20.0 °C) row is no longer needed - remove it.
$ws.Rows.Item(2).Delete()

# Add a new "Beschreibung" (description) column describing each phase of
# the remaining temperature profile points.
$ws.Range("C1").Value = "Beschreibung"
$ws.Range("C2").Value = "Aufheizen 1"
$ws.Range("C3").Value = "Aufheizen 2"
$ws.Range("C4").Value = "Aufheizen 3"
$ws.Range("C5").Value = "Halten"
$ws.Range("C6").Value = "Abkühlen 1"
$ws.Range("C7").Value = "Abkühlen 2"

# Match the formatting used by the existing columns (reuse the same style).
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("C2:C7").PasteSpecial(-4122)

$excel.CutCopyMode = 0
